$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3-16 down to 30-43, opening a 27-row gap (rows 3-29)
$ws.Rows("3:29").Insert()

# The inserted rows inherit formatting from the row above (the hyperlink row);
# fully clear that block so new cells fall back to the plain column style.
$ws.Range("A3:C29").Clear()

# --- New quiz content (pg 8 / pg9 sections) ---
$ws.Range("A4").Value = "pg 8"

$ws.Range("A5").Value = "Quid est Musica?"
$ws.Range("C5").Value = "What is music?"

$ws.Range("A6").Value = "Est bene canendi scientia."
$ws.Range("B6").Value = "Sie iſt ein kunſtrecht vnd wohl zu ſingen"
$ws.Range("C6").Value = "It is a well sung science / It is a good thing to sing."

$ws.Range("A7").Value = "Quotuplex est?"
$ws.Range("B7").Value = "Wievilfaltig iſt ſie"
$ws.Range("C7").Value = "How many is it? / How diverse is it? / How varied is it?"

$ws.Range("A8").Value = "Duplex -> Choralis est Figuralis"
$ws.Range("B8").Value = "Zwifaltig. Choralis oder gleichförmig/ und Figuralis oder vilförmlich."
$ws.Range("C8").Value = "Double -> Choral and Figured"

$ws.Range("A9").Value = "Quid est Musica Choralis?"
$ws.Range("C9").Value = "What is Choral Music"

$ws.Range("A10").Value = "Quae simplicem et uniformem in suis notulis servat mensuram."
$ws.Range("B10").Value = "Welche ein einfachen vnnd gleichförmigen Tact oder Menſur in jren Noten haltet."
$ws.Range("C10").Value = "That which is simple and uniform, keeps notes in their measure | Simple and uniform tact or mensur in your notes"

$ws.Range("A11").Value = "Quid est Musica Figuralis?"
$ws.Range("C11").Value = "What is musica figuralis?"

$ws.Range("A12").Value = "Quae diversam figurarum quantitatem, juxta varia praescripta signa exprimit."
$ws.Range("B12").Value = "Welche ein ungleiche gröſſe der Noten und Pauſen hat - nach mancherlen fürgeſetzten zeichen."
$ws.Range("C12").Value = "That which requires a different quantity of figures, according to the prescriptions of the standards of a variety of sounds. | That which has uneven sizes of notes and pauses"

$ws.Range("A13").Value = "Quot sunt praecipua capita, quibus tyro opus habet?"
$ws.Range("C13").Value = "What are the most important notes that the beginner needs?"

$ws.Range("A14").Value = " Quinque, Clavis, Vox, Cantus, Mutatio et Figura"
$ws.Range("C14").Value = "Five: Clavis (Music Key), voice, singing, mutation, and figure"

$ws.Range("A16").Value = "pg9"

$ws.Range("A18").Value = "Quid est clavis?"
$ws.Range("C18").Value = "What is the Clavis?"

$ws.Range("A19").Value = "Est vocis formandae index."
$ws.Range("C19").Value = "It is the voice to be formed index; It is the pointer to the voice that you want to sing. "

$ws.Range("A20").Value = "Quot sunt claves?"
$ws.Range("B20").Value = "wieviel sind musikschlüssel"
$ws.Range("C20").Value = "How are the musical keys?"

$ws.Range("A21").Value = "Viginti. Atque ex sequenti figura, quae vulgo Scala dicitur, patent."
$ws.Range("B21").Value = "Zweinzig und werden außnachfolgenden Täfelein ertände das gemeiniglich Scala (ein Leiter) genennet wird."
$ws.Range("C21").Value = "20. Now from the following figure, which is commonly known as the scale"

# Row 21 needs an explicit (non-auto) taller height to fit the German text
$ws.Rows("21:21").RowHeight = 54

# --- Update the view: scroll so row 10 is at the top, select C21 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C21").Select()
